# Master update 22 Nov 2020
# Re-orders a handful of Item Name / UOM pairs on Sheet1 (the BSL/BRAND/ISL rows
# stay put - only the Item Name (D) and UOM (E) text for a few rows is corrected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dinafex block (rows 3-5): 180mg/120mg/60mg -> 60mg/120mg/180mg
$ws.Range("D3").Value = "Dinafex 60mg Tablet"
$ws.Range("D5").Value = "Dinafex 180mg Tablet"

# Etorix block (rows 7-9): swap rows 8 and 9 (Item Name + UOM together)
$ws.Range("D8").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("E8").Value = "40's"
$ws.Range("D9").Value = "Etorix 90mg Tablet"
$ws.Range("E9").Value = "30's"

# Ketonic block (rows 14-16): rotate 14<-15, 15<-16, 16<-14
$ws.Range("D14").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E14").Value = "4's"
$ws.Range("D15").Value = "Ketonic 10mg Tablet"
$ws.Range("E15").Value = "20's"
$ws.Range("D16").Value = "Ketonic 30mg Injection"
$ws.Range("E16").Value = "5 's"

# Kynol block (rows 17-18): swap
$ws.Range("D17").Value = "Kynol TR 200mg Capsule"
$ws.Range("E17").Value = "30 's"
$ws.Range("D18").Value = "Kynol TR 100mg Capsule"
$ws.Range("E18").Value = "50 's"

# Zithrox block (rows 25,27): swap
$ws.Range("D25").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E25").Value = "30ml"
$ws.Range("D27").Value = "Zithrox 15ml Suspension"
$ws.Range("E27").Value = "15 ml"
